$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first two data rows (original USN 1MS21cs098 and 1MS21cs099),
# shifting the remaining rows up.
$ws.Range("A2:D3").EntireRow.Delete()

# After the shift, the rows that used to hold USN 1MS21cs103 .. 1MS21cs110
# are now rows 5 through 12. Remove them, keeping only the three rows for
# 1MS21cs100, 1MS21cs101, 1MS21cs102 (now rows 2-4).
$ws.Range("A5:D12").EntireRow.Delete()

$ws.Range("A1:D4").Select()
